$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# Header cell, styled like the other header cells (bold/centered/bordered)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Per-row timestamps
$ws.Range("F2").Value = "2021-10-05 13:41:21.842884"
$ws.Range("F3").Value = "2021-10-05 13:41:21.842896"
$ws.Range("F4").Value = "2021-10-05 13:41:21.842900"
$ws.Range("F5").Value = "2021-10-05 13:41:21.842904"
$ws.Range("F6").Value = "2021-10-05 13:41:21.842907"
$ws.Range("F7").Value = "2021-10-05 13:41:21.842910"
$ws.Range("F8").Value = "2021-10-05 13:41:21.842913"
$ws.Range("F9").Value = "2021-10-05 13:41:21.842916"
$ws.Range("F10").Value = "2021-10-05 13:41:21.842919"
